$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 0.083686
$ws.Range("N2").Value = 0.251058
$ws.Range("O2").Value = 0.08174459316063268
$ws.Range("P2").Value = 0.08174459316063269
$ws.Range("Q2").Value = 0.005964078057333334
$ws.Range("R2").Value = 0.053676702516
$ws.Range("S2").Value = 0.08174459316063268
$ws.Range("T2").Value = 0.08174459316063269

# Row 3 updates
$ws.Range("O3").Value = 0.5964052409947874
$ws.Range("P3").Value = 0.5964052409947875
$ws.Range("S3").Value = 0.5964052409947874
$ws.Range("T3").Value = 0.5964052409947875

# Row 4 updates
$ws.Range("M4").Value = 0.329494
$ws.Range("N4").Value = 0.988482
$ws.Range("O4").Value = 0.3218501658445798
$ws.Range("P4").Value = 0.3218501658445799
$ws.Range("Q4").Value = 0.02348215872933334
$ws.Range("R4").Value = 0.211339428564
$ws.Range("S4").Value = 0.3218501658445798
$ws.Range("T4").Value = 0.3218501658445799
